$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.06582034044206071
$ws.Range("C2").Value = 0.2316970354585097
$ws.Range("D2").Value = 0.09010395069034252
$ws.Range("E2").Value = 0.3001732011528386
$ws.Range("F2").Value = 0.3039234302155817

$ws.Range("B3").Value = 0.04030545542864885
$ws.Range("C3").Value = 0.117135837796785
$ws.Range("D3").Value = 0.03941846264558947
$ws.Range("E3").Value = 0.1985408336982331
$ws.Range("F3").Value = 0.2049225677834833

$ws.Range("B4").Value = 0.04535421744014654
$ws.Range("C4").Value = 0.1574291855501446
$ws.Range("D4").Value = 0.05660956138547054
$ws.Range("E4").Value = 0.2379276389692264
$ws.Range("F4").Value = 0.2558575142829206

$ws.Range("B5").Value = -0.03121797845395477
$ws.Range("C5").Value = 0.03121797845395477
$ws.Range("D5").Value = 0.001408506301122686
$ws.Range("E5").Value = 0.03753007195733424
$ws.Range("F5").Value = 0.02945994305395383
